$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reorder rows 3-5 and 7 (row 6 "Distributed PV" stays unchanged)
# New order: Hydro, Biomass, Wind, Distributed PV (unchanged), Utility-scale Solar

$ws.Range("A3").Value = "Hydro"
$ws.Range("B3").Value = 1.248976476713115
$ws.Range("C3").Value = 1122.23
$ws.Range("D3").Value = 1224.0985

$ws.Range("A4").Value = "Biomass"
$ws.Range("B4").Value = 18.48297686997073
$ws.Range("C4").Value = 12.3745
$ws.Range("D4").Value = 40.562

$ws.Range("A5").Value = "Wind"
$ws.Range("B5").Value = 7.981007756468905
$ws.Range("C5").Value = 1715.527
$ws.Range("D5").Value = 2936.4945

$ws.Range("A7").Value = "Utility-scale Solar"
$ws.Range("B7").Value = 54.75281499684422
$ws.Range("C7").Value = 134.704
$ws.Range("D7").Value = 2863.1915
